$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6803
$ws.Range("K3").Value = 7030
$ws.Range("K4").Value = 1451
$ws.Range("K6").Value = 7677
$ws.Range("K7").Value = 23461
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 207
$ws.Range("K5").Value = 62
$ws.Range("K7").Value = 714
$ws.Range("K8").Value = 1540
$ws.Range("K9").Value = 106
$ws.Range("K13").Value = 32
$ws.Range("K16").Value = 58
$ws.Range("K18").Value = 155
$ws.Range("K20").Value = 566
$ws.Range("K22").Value = 76
$ws.Range("K23").Value = 232
$ws.Range("K27").Value = 221
$ws.Range("K29").Value = 1273
$ws.Range("K31").Value = 259
$ws.Range("K33").Value = 1011
$ws.Range("K34").Value = 135
$ws.Range("K42").Value = 866
$ws.Range("K48").Value = 300
$ws.Range("K49").Value = 129
$ws.Range("K51").Value = 292
$ws.Range("K52").Value = 619
$ws.Range("K57").Value = 89
$ws.Range("K60").Value = 135
$ws.Range("K62").Value = 8
$ws.Range("K63").Value = 65
$ws.Range("K67").Value = 912
$ws.Range("K73").Value = 210
$ws.Range("K76").Value = 317
$ws.Range("K77").Value = 160
$ws.Range("K78").Value = 272
$ws.Range("K79").Value = 584
$ws.Range("K83").Value = 502
$ws.Range("K84").Value = 188
$ws.Range("K85").Value = 1087
$ws.Range("K88").Value = 254
$ws.Range("K89").Value = 350
$ws.Range("K90").Value = 223
$ws.Range("K93").Value = 88
$ws.Range("K96").Value = 251
$ws.Range("K98").Value = 118
$ws.Range("K99").Value = 400
$ws.Range("K101").Value = 23461
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 77
$ws.Range("K3").Value = 51
$ws.Range("K6").Value = 104
$ws.Range("K7").Value = 251
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 238
$ws.Range("K4").Value = 25
$ws.Range("K6").Value = 193
$ws.Range("K7").Value = 714
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 111
$ws.Range("K7").Value = 350
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 356
$ws.Range("K3").Value = 379
$ws.Range("K6").Value = 265
$ws.Range("K7").Value = 1087
$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K4").Value = 35
$ws.Range("K7").Value = 619
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 430
$ws.Range("K6").Value = 508
$ws.Range("K7").Value = 1540
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K6").Value = 115
$ws.Range("K7").Value = 502
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 256
$ws.Range("K6").Value = 319
$ws.Range("K7").Value = 1011
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 106
$ws.Range("K3").Value = 165
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 400
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 66
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 259
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 333
$ws.Range("K6").Value = 262
$ws.Range("K7").Value = 912
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 77
$ws.Range("K7").Value = 188
$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 63
$ws.Range("K7").Value = 129
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 361
$ws.Range("K3").Value = 453
$ws.Range("K6").Value = 368
$ws.Range("K7").Value = 1273
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 72
$ws.Range("K7").Value = 300
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K4").Value = 22
$ws.Range("K6").Value = 160
$ws.Range("K7").Value = 317
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K6").Value = 321
$ws.Range("K7").Value = 866
$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("K3").Value = 11
$ws.Range("K6").Value = 32
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 70
$ws.Range("K6").Value = 92
$ws.Range("K7").Value = 272
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 232
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 196
$ws.Range("K7").Value = 584
$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K2").Value = 31
$ws.Range("K6").Value = 53
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 183
$ws.Range("K7").Value = 566
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 155
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 88
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 53
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 135
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 118
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 106
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 210
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 61
$ws.Range("K7").Value = 207
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K2").Value = 65
$ws.Range("K3").Value = 78
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 254
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 62
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 58
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 221
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 63
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 223
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 79
$ws.Range("K7").Value = 292
$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 89
$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K2").Value = 44
$ws.Range("K7").Value = 135
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 76
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 160
$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 58
$ws = $wb.Worksheets.Item("Museum Campus")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 8
